# Applies the "Add files via upload" revision to EASYCADCOMMANDS.xlsx:
#   - Dimentions sheet: remove the placeholder "DBK" / "[DBK]..xxxxxxx" row
#   - Layout sheet: remove the placeholder "VVV" / "[VVV]..xxxxx" row
#   - General sheet: rename "FFF" / "[FFF]..Filter Of Selected Elements"
#     command to "FFFX" / "[FFFX]..Filter Of Selected Elements"
#   - Rename the trailing "Sheet1" tab to "Activiation"
#   - Move the active-tab / selection state from the old "Sheet1" tab back
#     to the "General" tab (first tab), updating each touched sheet's
#     remembered selection along the way

$wb = $excel.ActiveWorkbook

# --- Dimentions: drop the "DBK" row (row 6) ---------------------------------
$wsDim = $wb.Worksheets.Item("Dimentions")
$wsDim.Rows("6:6").Delete() | Out-Null
$wsDim.Range("A6:XFD6").Select() | Out-Null

# --- Layout: drop the "VVV" row (row 2) -------------------------------------
$wsLayout = $wb.Worksheets.Item("Layout")
$wsLayout.Rows("2:2").Delete() | Out-Null
$wsLayout.Range("A2:XFD2").Select() | Out-Null

# --- General: FFF -> FFFX ----------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("General")
$wsGeneral.Range("A1").Value = "FFFX"
$wsGeneral.Range("B1").Value = "[FFFX]..Filter Of Selected Elements"

# --- Rename the last sheet from "Sheet1" to "Activiation" ------------------
$wsLast = $wb.Worksheets.Item("Sheet1")
$wsLast.Name = "Activiation"
$wsLast.Range("G8").Select() | Out-Null

# --- Move the active tab back to "General", updating its selection ---------
$wsGeneral.Range("N6").Select() | Out-Null
